$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H48").Value = 0.75

$ws.Range("E49").Value = "Done"
$ws.Range("F49").Value = "Valid"
$ws.Range("H49").Value = 0.75

$ws.Range("E50").Value = "Done"
$ws.Range("F50").Value = "Valid"
$ws.Range("H50").Value = 0.75

$ws.Range("H51").Select()
